$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-6 with corrected (error-fixed) figures
# Row 2
$ws.Range("D2").Value = 1388
$ws.Range("E2").Value = 95
$ws.Range("F2").Value = 91
$ws.Range("G2").Value = 98
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 45
$ws.Range("J2").Value = 24
$ws.Range("K2").Value = 1557
$ws.Range("L2").Value = 577
$ws.Range("M2").Value = 980
$ws.Range("N2").Value = 774
$ws.Range("O2").Value = 207
$ws.Range("P2").Value = 50
$ws.Range("Q2").Value = 113
$ws.Range("R2").Value = -185
$ws.Range("S2").Value = 130
$ws.Range("T2").Value = 99
$ws.Range("U2").Value = 13
$ws.Range("V2").Value = 338
$ws.Range("W2").Value = 6.84
$ws.Range("X2").Value = 5.01
$ws.Range("Y2").Value = 5.98
$ws.Range("Z2").Value = 4.86
$ws.Range("AA2").Value = 58.87
$ws.Range("AB2").Value = 1449.04
$ws.Range("AC2").Value = 453
$ws.Range("AD2").Value = 14.83
$ws.Range("AE2").Value = 7735
$ws.Range("AF2").Value = 0.87
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 1.49
$ws.Range("AI2").Value = 22.07
$ws.Range("AJ2").Value = 10000000

# Row 3
$ws.Range("D3").Value = 1455
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 80
$ws.Range("G3").Value = 82
$ws.Range("H3").Value = 67
$ws.Range("I3").Value = 37
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 1551
$ws.Range("L3").Value = 537
$ws.Range("M3").Value = 1014
$ws.Range("N3").Value = 801
$ws.Range("O3").Value = 213
$ws.Range("P3").Value = 50
$ws.Range("Q3").Value = 129
$ws.Range("R3").Value = -110
$ws.Range("S3").Value = -23
$ws.Range("T3").Value = 93
$ws.Range("U3").Value = 36
$ws.Range("V3").Value = 340
$ws.Range("W3").Value = 5.48
$ws.Range("X3").Value = 4.62
$ws.Range("Y3").Value = 4.7
$ws.Range("Z3").Value = 4.33
$ws.Range("AA3").Value = 53
$ws.Range("AB3").Value = 1501.02
$ws.Range("AC3").Value = 370
$ws.Range("AD3").Value = 39.46
$ws.Range("AE3").Value = 8009
$ws.Range("AF3").Value = 1.82
$ws.Range("AG3").Value = 130
$ws.Range("AH3").Value = 0.89
$ws.Range("AI3").Value = 35.14
$ws.Range("AJ3").Value = 10000000

# Row 4
$ws.Range("D4").Value = 1633
$ws.Range("E4").Value = 122
$ws.Range("F4").Value = 122
$ws.Range("G4").Value = 119
$ws.Range("H4").Value = 99
$ws.Range("I4").Value = 56
$ws.Range("J4").Value = 43
$ws.Range("K4").Value = 1765
$ws.Range("L4").Value = 667
$ws.Range("M4").Value = 1098
$ws.Range("N4").Value = 844
$ws.Range("O4").Value = 254
$ws.Range("P4").Value = 50
$ws.Range("Q4").Value = 114
$ws.Range("R4").Value = -214
$ws.Range("S4").Value = 64
$ws.Range("T4").Value = 97
$ws.Range("U4").Value = 18
$ws.Range("V4").Value = 420
$ws.Range("W4").Value = 7.49
$ws.Range("X4").Value = 6.07
$ws.Range("Y4").Value = 6.81
$ws.Range("Z4").Value = 5.97
$ws.Range("AA4").Value = 60.72
$ws.Range("AB4").Value = 1586.91
$ws.Range("AC4").Value = 560
$ws.Range("AD4").Value = 19.91
$ws.Range("AE4").Value = 8443
$ws.Range("AF4").Value = 1.32
$ws.Range("AG4").Value = 110
$ws.Range("AH4").Value = 0.99
$ws.Range("AI4").Value = 19.64
$ws.Range("AJ4").Value = 10000000

# Row 5
$ws.Range("D5").Value = 1766
$ws.Range("E5").Value = 110
$ws.Range("F5").Value = 110
$ws.Range("G5").Value = 97
$ws.Range("H5").Value = 59
$ws.Range("I5").Value = 34
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 1846
$ws.Range("L5").Value = 680
$ws.Range("M5").Value = 1166
$ws.Range("N5").Value = 877
$ws.Range("O5").Value = 289
$ws.Range("P5").Value = 50
$ws.Range("Q5").Value = 59
$ws.Range("R5").Value = -95
$ws.Range("S5").Value = 26
$ws.Range("T5").Value = 97
$ws.Range("U5").Value = -38
$ws.Range("V5").Value = 399
$ws.Range("W5").Value = 6.23
$ws.Range("X5").Value = 3.32
$ws.Range("Y5").Value = 3.92
$ws.Range("Z5").Value = 3.25
$ws.Range("AA5").Value = 58.29
$ws.Range("AB5").Value = 1652.42
$ws.Range("AC5").Value = 338
$ws.Range("AD5").Value = 35.39
$ws.Range("AE5").Value = 8770
$ws.Range("AF5").Value = 1.36
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 1.26
$ws.Range("AI5").Value = 44.43
$ws.Range("AJ5").Value = 10000000

# Row 6
$ws.Range("D6").Value = 1861
$ws.Range("E6").Value = 33
$ws.Range("F6").Value = 33
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 3
$ws.Range("K6").Value = 1907
$ws.Range("L6").Value = 748
$ws.Range("M6").Value = 1159
$ws.Range("N6").Value = 869
$ws.Range("P6").Value = 50
$ws.Range("Q6").Value = 98
$ws.Range("R6").Value = -96
$ws.Range("S6").Value = 30
$ws.Range("T6").Value = 81
$ws.Range("U6").Value = 17
$ws.Range("V6").Value = 449
$ws.Range("W6").Value = 1.78
$ws.Range("X6").Value = 0.25
$ws.Range("Y6").Value = 0.38
$ws.Range("Z6").Value = 0.24
$ws.Range("AA6").Value = 64.55
$ws.Range("AB6").Value = 1631.47
$ws.Range("AC6").Value = 33
$ws.Range("AD6").Value = 338.55
$ws.Range("AE6").Value = 8686
$ws.Range("AF6").Value = 1.28
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 0.9
$ws.Range("AI6").Value = 305
$ws.Range("AJ6").Value = 10000000

# Rows 7-9 no longer have reliable data for columns D:AI -> clear them
$ws.Range("D7:AI9").ClearContents()
